$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: title + link
$ws.Range("D9").Value = "맞는 시험 vs. 틀린 시험"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/right-exam-wrong-exam/#utm_source=rss&utm_medium=rss&utm_campaign=right-exam-wrong-exam"

# Row 26: title only
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 51: title + link
$ws.Range("D51").Value = "마케팅 용어, 앰부시 의미"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EB%A7%88%EC%BC%80%ED%8C%85-%EC%9A%A9%EC%96%B4-%EC%95%B0%EB%B6%80%EC%8B%9C-%EC%9D%98%EB%AF%B8"

# Row 52: title only
$ws.Range("D52").Value = "Relative Risk Regression"
